$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-ambiguous updates (Coin names, Links, thousands-dotted prices,
# and Volume(1h) percentages) - assigning .Value keeps these as text since Excel
# cannot parse them as a single number.
$ws.Range('D2').Value = '96.470.65'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '3.712.11'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('E6').Value = '  +7.57%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D11').Value = '3.711.93'
$ws.Range('E11').Value = '  +3.32%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  +5.79%  '
$ws.Range('D15').Value = '4.404.86'
$ws.Range('E15').Value = '  +3.31%  '
$ws.Range('E16').Value = '  +2.86%  '
$ws.Range('D17').Value = '96.227.98'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('E18').Value = '  +15.64%  '
$ws.Range('D19').Value = '3.714.29'
$ws.Range('E19').Value = '  +3.71%  '
$ws.Range('E20').Value = '  +4.29%  '
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  -3.33%  '
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('E24').Value = '  -0.26%  '
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  -1.21%  '
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('E28').Value = '  +0.82%  '
$ws.Range('E29').Value = '  -6.93%  '
$ws.Range('E30').Value = '  +3.23%  '
$ws.Range('E31').Value = '  +2.11%  '
$ws.Range('E32').Value = '  +0.13%  '
$ws.Range('E33').Value = '  +10.03%  '
$ws.Range('E34').Value = '  -1.72%  '
$ws.Range('E35').Value = '  +8.31%  '
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('E38').Value = '  +2.45%  '
$ws.Range('E39').Value = '  +1.12%  '
$ws.Range('E40').Value = '  +16.10%  '
$ws.Range('E41').Value = '  +4.48%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E42').Value = '  +5.05%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E43').Value = '  +21.92%  '
$ws.Range('E44').Value = '  +2.26%  '
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('E47').Value = '  -2.20%  '
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('E51').Value = '  +2.40%  '

# Price cells whose new text is a single valid decimal number (e.g. "238.08", "1.00").
# A plain .Value assignment would let Excel auto-convert these to a numeric cell
# (changing the cell type and, for "clean" values like "1.00"/"45.00", silently
# dropping the trailing zeros). Force literal text via a leading apostrophe - Excel's
# native quote-prefix input - via .Formula so the stored value is text, matching source.
$ws.Range('D5').Formula = "'238.08"
$ws.Range('D6').Formula = "'1.94"
$ws.Range('D7').Formula = "'654.42"
$ws.Range('D8').Formula = "'0.422"
$ws.Range('D10').Formula = "'1.00"
$ws.Range('D12').Formula = "'45.00"
$ws.Range('D14').Formula = "'6.83"
$ws.Range('D18').Formula = "'8.99"
$ws.Range('D20').Formula = "'19.07"
$ws.Range('D21').Formula = "'12.80"
$ws.Range('D22').Formula = "'0.530"
$ws.Range('D23').Formula = "'525.47"
$ws.Range('D24').Formula = "'3.49"
$ws.Range('D27').Formula = "'102.63"
$ws.Range('D28').Formula = "'13.43"
$ws.Range('D30').Formula = "'12.48"
$ws.Range('D31').Formula = "'3.06"
$ws.Range('D33').Formula = "'1.90"
$ws.Range('D34').Formula = "'0.187"
$ws.Range('D35').Formula = "'669.95"
$ws.Range('D36').Formula = "'32.76"
$ws.Range('D39').Formula = "'8.87"
$ws.Range('D40').Formula = "'7.13"
$ws.Range('D42').Formula = "'0.977"
$ws.Range('D43').Formula = "'39.88"
$ws.Range('D46').Formula = "'0.0460"
$ws.Range('D47').Formula = "'0.439"
$ws.Range('D49').Formula = "'23.61"
$ws.Range('D50').Formula = "'8.61"
